# Updates the cryptos list with refreshed price / 1h-volume figures
# (and fixes the coin ordering for three swapped pairs of rows), matching
# the "Updated cryptos list ... with GitHub Actions" commit.
#
# All values are plain text in this sheet (prices such as "43.411.96" use
# dots as thousands separators, not decimals, and volume change values keep
# their padding spaces), so every write is prefixed with a leading
# apostrophe to force Excel to store it as text instead of silently
# re-parsing it as a number (which would mangle things like "1.00" -> 1,
# "6.60" -> 6.6, or "43.411.96" -> a date/garbage).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.411.96"
$ws.Range("E2").Value = "'  -1.38%  "
$ws.Range("D3").Value = "'2.365.58"
$ws.Range("E3").Value = "'  +5.31%  "
$ws.Range("E4").Value = "'  +0.07%  "
$ws.Range("D5").Value = "'233.32"
$ws.Range("E5").Value = "'  +1.02%  "
$ws.Range("D6").Value = "'0.651"
$ws.Range("E6").Value = "'  +1.36%  "
$ws.Range("D7").Value = "'68.92"
$ws.Range("E7").Value = "'  +9.06%  "
$ws.Range("E8").Value = "'  -0.05%  "
$ws.Range("D9").Value = "'0.459"
$ws.Range("E9").Value = "'  +2.34%  "
$ws.Range("D10").Value = "'0.0969"
$ws.Range("E10").Value = "'  -1.00%  "
$ws.Range("D11").Value = "'57.04"
$ws.Range("E11").Value = "'  -0.47%  "
$ws.Range("D12").Value = "'26.34"
$ws.Range("E12").Value = "'  -0.21%  "
$ws.Range("D13").Value = "'2.721.32"
$ws.Range("E13").Value = "'  +5.62%  "
$ws.Range("D14").Value = "'0.105"
$ws.Range("E14").Value = "'  -0.46%  "
$ws.Range("E16").Value = "'  +1.49%  "
$ws.Range("D17").Value = "'0.851"
$ws.Range("E17").Value = "'  +2.98%  "
$ws.Range("D18").Value = "'2.364.21"
$ws.Range("E18").Value = "'  +5.43%  "
$ws.Range("D19").Value = "'43.413.67"
$ws.Range("E19").Value = "'  -0.95%  "
$ws.Range("D20").Value = "'0.0₃0985"
$ws.Range("E20").Value = "'  -0.39%  "
$ws.Range("D21").Value = "'6.35"
$ws.Range("E21").Value = "'  +4.73%  "
$ws.Range("D22").Value = "'73.75"
$ws.Range("E22").Value = "'  +1.56%  "
$ws.Range("D23").Value = "'250.66"
$ws.Range("E23").Value = "'  +1.28%  "
$ws.Range("D24").Value = "'3.88"
$ws.Range("E24").Value = "'  +16.38%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "'  -0.15%  "
$ws.Range("E26").Value = "'  +1.94%  "
$ws.Range("D27").Value = "'2.28"
$ws.Range("E27").Value = "'  +1.78%  "
$ws.Range("D28").Value = "'22.74"
$ws.Range("E28").Value = "'  +8.55%  "
$ws.Range("D29").Value = "'9.94"
$ws.Range("E29").Value = "'  +1.13%  "
$ws.Range("D30").Value = "'172.43"
$ws.Range("E30").Value = "'  +0.47%  "
$ws.Range("D31").Value = "'1.56"
$ws.Range("E31").Value = "'  +9.96%  "
$ws.Range("E32").Value = "'  -7.81%  "
$ws.Range("E33").Value = "'  +1.39%  "
$ws.Range("E34").Value = "'  +4.47%  "
$ws.Range("D35").Value = "'0.0689"
$ws.Range("E35").Value = "'  +0.31%  "
$ws.Range("D36").Value = "'5.07"
$ws.Range("E36").Value = "'  +2.56%  "
$ws.Range("B37").Value = "'THORChain"
$ws.Range("C37").Value = "'https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D37").Value = "'6.60"
$ws.Range("E37").Value = "'  +2.82%  "
$ws.Range("B38").Value = "'LidoDAOToken"
$ws.Range("C38").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "'2.45"
$ws.Range("E38").Value = "'  +7.55%  "
$ws.Range("D39").Value = "'3.61"
$ws.Range("E39").Value = "'  -1.29%  "
$ws.Range("E40").Value = "'  -0.23%  "
$ws.Range("E41").Value = "'  +8.09%  "
$ws.Range("E42").Value = "'  +0.08%  "
$ws.Range("D43").Value = "'18.41"
$ws.Range("E43").Value = "'  +7.78%  "
$ws.Range("D44").Value = "'1.18"
$ws.Range("E44").Value = "'  +10.19%  "
$ws.Range("B45").Value = "'TrustWalletToken"
$ws.Range("C45").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").Value = "'1.22"
$ws.Range("E45").Value = "'  +2.40%  "
$ws.Range("B46").Value = "'Aave"
$ws.Range("C46").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'98.93"
$ws.Range("E46").Value = "'  +1.88%  "
$ws.Range("D47").Value = "'4.46"
$ws.Range("E47").Value = "'  +2.95%  "
$ws.Range("D48").Value = "'0.0955"
$ws.Range("E48").Value = "'  +1.56%  "
$ws.Range("D49").Value = "'1.448.53"
$ws.Range("E49").Value = "'  +1.36%  "
$ws.Range("D50").Value = "'2.588.43"
$ws.Range("E50").Value = "'  +5.56%  "
$ws.Range("B51").Value = "'Celestia"
$ws.Range("C51").Value = "'https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D51").Value = "'9.90"
$ws.Range("E51").Value = "'  +0.92%  "
